$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Finance department -> DEP001 / division finance
$ws.Range("C2").Value = "DEP001"
$ws.Range("E2").Value = "division finance"

# Row 3: IT department -> branch_code B001, code DEP002 / division it
$ws.Range("B3").Value = "B001"
$ws.Range("C3").Value = "DEP002"
$ws.Range("E3").Value = "division it"
